# Adds R and immediate (I-TYPE) control-signal data to rows 38-45 of the
# decoding spreadsheet (columns AW:BB), mirroring the values already present
# in row 37 (ADDI), and updates the sheet's frozen-pane/selection view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 38; $r -le 45; $r++) {
    $ws.Range("AW$r").Value = 1
    $ws.Range("AX$r").Value = 1
    $ws.Range("AY$r").Value = 1
    $ws.Range("AZ$r").Value = 0
    $ws.Range("BA$r").Value = "00"
    $ws.Range("BB$r").Value = "N"
}

# Update the view: move the active selection to the bottom-right pane's new
# location (scrolled further down the sheet).
$ws.Range("BB42").Select()
